$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The report is organised in blocks of 10 rows (job_id 71,72,74,75,76,77,78,
# 79,80,81). A new job "82" row (with its own "vecinos" count) is inserted
# right after each block, pushing the rest of the sheet down by one row.
# Working from the bottom of the sheet upward so the earlier row numbers
# below each insertion point stay valid while we go.
$inserts = @(
    @{ Row = 52; Vecinos = "55" },
    @{ Row = 42; Vecinos = "50" },
    @{ Row = 32; Vecinos = "61" },
    @{ Row = 22; Vecinos = "78" },
    @{ Row = 12; Vecinos = "63" }
)

foreach ($ins in $inserts) {
    $r = $ins.Row
    $ws.Range("A$r").EntireRow.Insert()
    $ws.Range("A$r").Value = 82

    # Write the vecinos count as genuine text (matching the rest of column B)
    # without perturbing any cell styles: stage it as a TEXT() formula, then
    # collapse it down to a literal value via a values-only paste.
    $cell = $ws.Range("B$r")
    $cell.Formula = '=TEXT(' + $ins.Vecinos + ',"0")'
    $cell.Copy()
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $excel.CutCopyMode = 0
}
